$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.080.52'
$ws.Range('E2').Value = '  +0.91%  '

$ws.Range('D3').Value = '3.006.87'
$ws.Range('E3').Value = '  +2.97%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '353.69'
$ws.Range('E5').Value = '  -0.14%  '

$ws.Range('D6').Value = '106.98'
$ws.Range('E6').Value = '  -2.47%  '

$ws.Range('D7').Value = '0.557'
$ws.Range('E7').Value = '  -0.72%  '

$ws.Range('E8').Value = '  +0.25%  '

$ws.Range('D9').Value = '0.610'
$ws.Range('E9').Value = '  -3.45%  '

$ws.Range('D10').Value = '38.07'
$ws.Range('E10').Value = '  -2.61%  '

$ws.Range('E11').Value = '  +2.25%  '

$ws.Range('D12').Value = '0.0854'
$ws.Range('E12').Value = '  -4.36%  '

$ws.Range('D13').Value = '18.98'
$ws.Range('E13').Value = '  -3.77%  '

$ws.Range('D14').Value = '3.478.14'
$ws.Range('E14').Value = '  +2.90%  '

$ws.Range('E15').Value = '  -3.70%  '

$ws.Range('D16').Value = '3.002.68'
$ws.Range('E16').Value = '  +2.51%  '

$ws.Range('E17').Value = '  +3.30%  '

$ws.Range('D18').Value = '52.110.07'
$ws.Range('E18').Value = '  +0.87%  '

$ws.Range('D19').Value = '3.43'
$ws.Range('E19').Value = '  +5.17%  '

$ws.Range('D20').Value = '7.47'
$ws.Range('E20').Value = '  -1.46%  '

$ws.Range('E21').Value = '  -4.55%  '

$ws.Range('D22').Value = '0.0₃0970'
$ws.Range('E22').Value = '  -1.26%  '

$ws.Range('D23').Value = '68.99'
$ws.Range('E23').Value = '  -2.58%  '

$ws.Range('D24').Value = '264.28'
$ws.Range('E24').Value = '  -2.04%  '

$ws.Range('E25').Value = '  -3.26%  '

$ws.Range('D26').Value = '0.177'
$ws.Range('E26').Value = '  -2.98%  '

$ws.Range('D27').Value = '26.93'
$ws.Range('E27').Value = '  -0.86%  '

$ws.Range('E28').Value = '  -0.05%  '

$ws.Range('D29').Value = '7.39'
$ws.Range('E29').Value = '  -2.13%  '

$ws.Range('E30').Value = '  -1.18%  '

$ws.Range('D31').Value = '6.36'
$ws.Range('E31').Value = '  +3.54%  '

$ws.Range('D32').Value = '10.19'
$ws.Range('E32').Value = '  -3.90%  '

$ws.Range('D33').Value = '35.96'
$ws.Range('E33').Value = '  -5.56%  '

$ws.Range('E34').Value = '  +8.69%  '

$ws.Range('D35').Value = '51.09'
$ws.Range('E35').Value = '  -2.26%  '

$ws.Range('D36').Value = '0.0438'
$ws.Range('E36').Value = '  -0.33%  '

$ws.Range('E37').Value = '  -0.08%  '

$ws.Range('D38').Value = '3.35'
$ws.Range('E38').Value = '  +3.33%  '

$ws.Range('D39').Value = '2.83'
$ws.Range('E39').Value = '  +4.00%  '

$ws.Range('D40').Value = '1.96'
$ws.Range('E40').Value = '  -2.51%  '

$ws.Range('D41').Value = '17.52'
$ws.Range('E41').Value = '  -4.12%  '

$ws.Range('E42').Value = '  -1.19%  '

$ws.Range('D43').Value = '124.50'
$ws.Range('E43').Value = '  +6.80%  '

$ws.Range('D44').Value = '22.80'
$ws.Range('E44').Value = '  -1.43%  '

$ws.Range('D45').Value = '2.12'
$ws.Range('E45').Value = '  -2.27%  '

$ws.Range('D46').Value = '2.121.01'
$ws.Range('E46').Value = '  -0.89%  '

$ws.Range('E47').Value = '  -3.76%  '

$ws.Range('E48').Value = '  -6.06%  '

$ws.Range('D49').Value = '3.300.89'
$ws.Range('E49').Value = '  +2.82%  '

$ws.Range('D50').Value = '0.242'
$ws.Range('E50').Value = '  -3.32%  '

$ws.Range('D51').Value = '0.0333'
$ws.Range('E51').Value = '  +1.04%  '
